$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal TEXT value even when it looks like a number,
# mirroring the source data (all D/E values in this sheet are authored as text).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    if ($text -match "^[+-]?[0-9]*\.?[0-9]+$") {
        $rng.NumberFormat = "@"
        $rng.Value = $text
        $rng.ClearFormats()
    } else {
        $rng.Value = $text
    }
}

Set-TextValue 'D2' '66.777.81'
Set-TextValue 'E2' '  +0.77%  '
Set-TextValue 'D3' '3.495.82'
Set-TextValue 'E3' '  +0.15%  '
Set-TextValue 'E4' '  -0.02%  '
Set-TextValue 'D5' '594.63'
Set-TextValue 'E5' '  -0.15%  '
Set-TextValue 'D6' '169.43'
Set-TextValue 'E6' '  -0.02%  '
Set-TextValue 'E7' '  +0.01%  '
Set-TextValue 'E8' '  +1.90%  '
Set-TextValue 'E9' '  +7.17%  '
Set-TextValue 'D10' '7.33'
Set-TextValue 'E10' '  +0.62%  '
Set-TextValue 'E11' '  -0.54%  '
Set-TextValue 'D12' '4.106.26'
Set-TextValue 'E12' '  +0.23%  '
Set-TextValue 'E13' '  -0.35%  '
Set-TextValue 'D14' '28.24'
Set-TextValue 'E14' '  +1.35%  '
Set-TextValue 'E15' '  +2.59%  '
Set-TextValue 'D16' '66.793.45'
Set-TextValue 'E16' '  +0.84%  '
Set-TextValue 'D17' '3.498.95'
Set-TextValue 'E17' '  -0.13%  '
Set-TextValue 'E18' '  +0.56%  '
Set-TextValue 'D19' '14.07'
Set-TextValue 'E19' '  +0.17%  '
Set-TextValue 'D20' '394.70'
Set-TextValue 'E20' '  +1.87%  '
Set-TextValue 'D21' '7.94'
Set-TextValue 'E21' '  -0.69%  '
Set-TextValue 'D22' '73.43'
Set-TextValue 'E22' '  +0.65%  '
Set-TextValue 'D24' '0.536'
Set-TextValue 'E24' '  +1.70%  '
Set-TextValue 'E25' '  +1.13%  '
Set-TextValue 'E26' '  +0.82%  '
Set-TextValue 'E27' '  +0.17%  '
Set-TextValue 'D28' '1.00'
Set-TextValue 'E28' '  -0.07%  '
Set-TextValue 'D29' '6.29'
Set-TextValue 'E29' '  -1.11%  '
Set-TextValue 'D30' '1.45'
Set-TextValue 'E30' '  -0.96%  '
Set-TextValue 'D31' '2.07'
Set-TextValue 'E31' '  -0.24%  '
Set-TextValue 'D32' '23.99'
Set-TextValue 'E32' '  +2.35%  '
Set-TextValue 'D33' '7.37'
Set-TextValue 'E33' '  -0.69%  '
Set-TextValue 'E34' '  +3.76%  '
Set-TextValue 'D35' '163.34'
Set-TextValue 'E35' '  +1.73%  '
Set-TextValue 'D36' '0.892'
Set-TextValue 'E36' '  -1.23%  '
Set-TextValue 'E37' '  -0.56%  '
Set-TextValue 'D38' '6.85'
Set-TextValue 'E38' '  +3.11%  '
Set-TextValue 'D39' '4.71'
Set-TextValue 'E39' '  +3.36%  '
Set-TextValue 'E40' '  -0.72%  '
Set-TextValue 'D41' '26.31'
Set-TextValue 'E41' '  -0.55%  '
Set-TextValue 'D42' '2.824.43'
Set-TextValue 'E42' '  +0.72%  '

# Row 43 <-> Row 44: coin identities swapped (ranking reorder) plus refreshed price/volume data
Set-TextValue 'B43' 'dogwifhat'
Set-TextValue 'C43' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D43' '2.62'
Set-TextValue 'E43' '  +5.31%  '
Set-TextValue 'B44' 'InjectiveProtocol'
Set-TextValue 'C44' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D44' '26.86'
Set-TextValue 'E44' '  -1.39%  '

Set-TextValue 'D45' '42.73'
Set-TextValue 'E45' '  -1.53%  '
Set-TextValue 'E46' '  -0.47%  '
Set-TextValue 'D47' '342.17'
Set-TextValue 'E48' '  +1.26%  '
Set-TextValue 'D49' '33.64'
Set-TextValue 'E49' '  +3.22%  '
Set-TextValue 'D50' '0.853'
Set-TextValue 'E50' '  +0.36%  '
Set-TextValue 'D51' '6.50'
Set-TextValue 'E51' '  +0.87%  '
